$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

$statQuery = @'
MATCH (p:program)<--(s:study)<--(c:case)<--(demo:demographic), (c)<--(diag:diagnosis), (c)<--(r:registration)
MATCH (samp:sample)-->(c)
WHERE s.clinical_study_designation IN ['MGT01'] and samp.sample_site in['Mammary Gland']
OPTIONAL MATCH (cf:file)-[*]->(c)
OPTIONAL MATCH (sf:file)-->(s)
RETURN
	count(distinct p) AS Programs,
    count(distinct s) AS Studies,
    count(distinct c) AS Cases,
    count(distinct samp) AS Samples,
    count(distinct cf) AS `Case Files`,
    count(distinct sf) AS `Study Files`

    
'@

$neo4jName = @'
TC44_Canine_Study_MGT01_SampleSite_MammaryGland_Neo4jData.xlsx
'@

$webName = @'
TC44_Canine_Study_MGT01_SampleSite_MammaryGland_WebData.xlsx
'@

$filesTabQuery = @'
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
MATCH (r:registration)-->(c)
MATCH (f)-[*]->(samp:sample)
WHERE s.clinical_study_designation IN ['MGT01'] and samp.sample_site in['Mammary Gland']
OPTIONAL MATCH (f)-[*]->(samp:sample)
WITH
        DISTINCT f, parent, c, demo, diag, s, samp,
        ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
        toInteger(floor(log(f.file_size)/log(1024))) as i,
        2 as precision
WITH
        f, parent, c, demo, diag, s, samp,
        f.file_size /(1024^i) AS value,
        10^precision AS factor,
        units[i] as unit
WITH
        f, parent, c, demo, diag, s, samp, unit,
        round(factor * value)/factor AS size
RETURN
        coalesce(f.file_name, '') AS `File Name`,
        coalesce(f.file_format, '') AS `Format`,
        coalesce(f.file_type, '') AS `File Type`,
        CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
        coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(samp.sample_id, '') AS `Sample ID`,
        coalesce(c.case_id, '') AS `Case ID`,
        coalesce(demo.breed,'') AS Breed ,
        coalesce(diag.disease_term,'') AS Diagnosis
        order by f.file_name asc
        limit 100
'@

# Row 2 (CasesTab): StatQuery (C2) unchanged, update dbExcel/WebExcel file names
$ws.Range("D2").Value = $neo4jName
$ws.Range("E2").Value = $webName

# Row 3 (SamplesTab): update StatQuery + file names
$ws.Range("C3").Value = $statQuery
$ws.Range("D3").Value = $neo4jName
$ws.Range("E3").Value = $webName

# Row 4 (FilesTab): update query (B4), StatQuery (C4) + file names
$ws.Range("B4").Value = $filesTabQuery
$ws.Range("C4").Value = $statQuery
$ws.Range("D4").Value = $neo4jName
$ws.Range("E4").Value = $webName

# Row 5 (StudyFilesTab): update StatQuery + file names
$ws.Range("C5").Value = $statQuery
$ws.Range("D5").Value = $neo4jName
$ws.Range("E5").Value = $webName

# Update selected cell to D5
$ws.Range("D5").Select() | Out-Null
